# Gestion des délégataires par communes
$wb = $excel.ActiveWorkbook

$wsCommunes = $wb.Worksheets.Item("Financements")
$wsNotice = $wb.Worksheets.Item("Notice")

# Update the last-selected cell on the "Notice" sheet (B10 -> B39), then
# switch back to the main sheet so it stays the active tab.
$wsNotice.Range("B39").Select()
$wsCommunes.Activate()

# Rename the first sheet from "Financements" to "Communes"
$wsCommunes.Name = "Communes"
